# Applies the "cryptos list" price/volume refresh described in the commit
# message ("Updated cryptos list on Wed May  8 21:39:22 UTC 2024 with GitHub
# Actions") by writing the new Price (D) / Volume(1h) (E) values, and for the
# two rows whose coins swapped rank (39/40) the Coin (B) and Link (C) values
# too.
#
# All of the source cells are plain text (t="inlineStr"), including values
# that look numeric ("1.00", "587.45", ...). A plain `Range.Value = "1.00"`
# assignment would get auto-coerced by Excel into the *number* 1, which is not
# what the diff wants. To keep every written cell as literal text - without
# leaving behind a new NumberFormat/style the way a quote-prefixed or
# Text-formatted cell would - each value is entered as a `="..."` text formula
# and then immediately collapsed back down to a plain value via Copy +
# PasteSpecial(xlPasteValues).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        [string]$Addr,
        [string]$Text
    )
    $cell = $ws.Range($Addr)
    $escaped = $Text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 2: D2, E2
Set-CellText 'D2' '61.604.27'
Set-CellText 'E2' '  -2.23%  '

# Row 3: D3, E3
Set-CellText 'D3' '2.971.96'
Set-CellText 'E3' '  -2.63%  '

# Row 4: D4, E4
Set-CellText 'D4' '1.00'
Set-CellText 'E4' '  -0.02%  '

# Row 5: D5, E5
Set-CellText 'D5' '587.45'
Set-CellText 'E5' '  +0.79%  '

# Row 6: D6, E6
Set-CellText 'D6' '141.37'
Set-CellText 'E6' '  -6.63%  '

# Row 7: D7, E7
Set-CellText 'D7' '1.00'
Set-CellText 'E7' '  +0.01%  '

# Row 8: D8, E8
Set-CellText 'D8' '0.519'
Set-CellText 'E8' '  -2.85%  '

# Row 9: D9, E9
Set-CellText 'D9' '2.968.64'
Set-CellText 'E9' '  -2.61%  '

# Row 10: E10
Set-CellText 'E10' '  -6.24%  '

# Row 11: D11, E11
Set-CellText 'D11' '5.75'
Set-CellText 'E11' '  -1.44%  '

# Row 12: E12
Set-CellText 'E12' '  +2.24%  '

# Row 13: E13
Set-CellText 'E13' '  -3.66%  '

# Row 14: D14, E14
Set-CellText 'D14' '33.95'
Set-CellText 'E14' '  -6.05%  '

# Row 15: E15
Set-CellText 'E15' '  +1.37%  '

# Row 16: D16, E16
Set-CellText 'D16' '3.460.62'
Set-CellText 'E16' '  -2.63%  '

# Row 17: D17, E17
Set-CellText 'D17' '6.99'
Set-CellText 'E17' '  -1.94%  '

# Row 18: D18, E18
Set-CellText 'D18' '61.595.35'
Set-CellText 'E18' '  -2.25%  '

# Row 19: D19, E19
Set-CellText 'D19' '2.972.02'
Set-CellText 'E19' '  -2.75%  '

# Row 20: D20, E20
Set-CellText 'D20' '449.76'
Set-CellText 'E20' '  -6.40%  '

# Row 21: D21, E21
Set-CellText 'D21' '13.85'
Set-CellText 'E21' '  -3.12%  '

# Row 22: E22
Set-CellText 'E22' '  -3.57%  '

# Row 23: E23
Set-CellText 'E23' '  -2.71%  '

# Row 24: D24, E24
Set-CellText 'D24' '81.04'
Set-CellText 'E24' '  -1.04%  '

# Row 25: E25
Set-CellText 'E25' '  -4.30%  '

# Row 26: E26
Set-CellText 'E26' '  -10.49%  '

# Row 28: E28
Set-CellText 'E28' '  -7.78%  '

# Row 29: E29
Set-CellText 'E29' '  -0.07%  '

# Row 30: D30, E30
Set-CellText 'D30' '2.63'
Set-CellText 'E30' '  -1.15%  '

# Row 31: D31, E31
Set-CellText 'D31' '6.84'
Set-CellText 'E31' '  -7.34%  '

# Row 32: E32
Set-CellText 'E32' '  -6.57%  '

# Row 33: D33, E33
Set-CellText 'D33' '27.15'
Set-CellText 'E33' '  -2.30%  '

# Row 34: E34
Set-CellText 'E34' '  -3.56%  '

# Row 35: E35
Set-CellText 'E35' '  -5.01%  '

# Row 36: D36, E36
Set-CellText 'D36' '0.0₃0772'
Set-CellText 'E36' '  -4.94%  '

# Row 37: E37
Set-CellText 'E37' '  -4.11%  '

# Row 38: E38
Set-CellText 'E38' '  -5.54%  '

# Row 39: B39, C39, D39, E39
Set-CellText 'B39' 'OKB'
Set-CellText 'C39' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText 'D39' '50.12'
Set-CellText 'E39' '  -0.54%  '

# Row 40: B40, C40, D40, E40
Set-CellText 'B40' 'Cosmos'
Set-CellText 'C40' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText 'D40' '9.13'
Set-CellText 'E40' '  -0.71%  '

# Row 41: D41, E41
Set-CellText 'D41' '0.118'
Set-CellText 'E41' '  +2.65%  '

# Row 42: D42, E42
Set-CellText 'D42' '2.79'
Set-CellText 'E42' '  -11.76%  '

# Row 43: D43, E43
Set-CellText 'D43' '386.99'
Set-CellText 'E43' '  -9.53%  '

# Row 44: D44, E44
Set-CellText 'D44' '0.0353'
Set-CellText 'E44' '  -2.40%  '

# Row 45: D45, E45
Set-CellText 'D45' '2.720.41'
Set-CellText 'E45' '  -4.36%  '

# Row 46: E46
Set-CellText 'E46' '  -8.61%  '

# Row 47: D47, E47
Set-CellText 'D47' '36.94'
Set-CellText 'E47' '  -2.70%  '

# Row 48: E48
Set-CellText 'E48' '  +2.13%  '

# Row 49: E49
Set-CellText 'E49' '  +0.07%  '

# Row 50: D50, E50
Set-CellText 'D50' '0.107'
Set-CellText 'E50' '  -1.85%  '

# Row 51: E51
Set-CellText 'E51' '  -1.64%  '

$excel.CutCopyMode = $false

Write-Output "Updated $(49) coin rows (83 cells) in Sheet1."
